$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1 & 2: two runs (“客户服务子系统可通过调用系统的功能与数据，为经纪人的客户提供服务。”
#        and “该系统可以”) pick up the explicit Helvetica/宋体 run formatting
#        (ascii/hAnsi/cs Helvetica, eastAsia 宋体, color 333333, kern 0) that the
#        rest of the document already uses, while keeping their existing
#        w:hint="eastAsia" + sz/szCs.
#
# NB: Range.Font.NameAscii/NameFarEast/NameOther/NameBi (when set directly on a
# Find-produced Range/Selection) mutate every run's rFonts in the whole
# enclosing paragraph in this host, not just the matched text. Driving the
# same assignment through Find.Replacement.Font instead (classic "formatted
# Find & Replace") applies it only to the matched span, splitting/merging
# runs exactly the way Word does.
# ---------------------------------------------------------------------------

function Set-HelveticaRun([string]$text) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Text = $text
    $find.Replacement.Font.NameAscii = "Helvetica"
    $find.Replacement.Font.NameFarEast = "宋体"
    $find.Replacement.Font.NameOther = "Helvetica"
    $find.Replacement.Font.NameBi = "Helvetica"
    $find.Replacement.Font.Color = 3355443
    $find.Replacement.Font.Kerning = 0
    $find.Replacement.Text = $text
    $find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, $text, 2) | Out-Null
}

Set-HelveticaRun("客户服务子系统可通过调用系统的功能与数据，为经纪人的客户提供服务。")
Set-HelveticaRun("该系统可以")

# ---------------------------------------------------------------------------
# 3 & 4: the editor's cursor ended up between “个股资讯” and
# “、市场公告和研究报告等信息内容” — Word records that as the (single,
# document-wide) "_GoBack" bookmark, splitting the run it lands in. Adding a
# bookmark named "_GoBack" here both creates the split/insertion in the diff
# and implicitly relocates the bookmark away from its old spot in the final
# (now fully empty) paragraph.
# ---------------------------------------------------------------------------

$splitPoint = $d.Content
$splitPoint.Find.Execute("提供信息服务，包括今日要闻、分类新闻、个股资讯", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPoint.Collapse(0)
$d.Bookmarks.Add("_GoBack", $splitPoint) | Out-Null
